$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.47'
$ws.Range("E2").Value = '''0.71%'
$ws.Range("G2").Value = '''23'
$ws.Range("E3").Value = '''6.96%'
$ws.Range("G3").Value = '''23'
$ws.Range("D4").Value = '''5.195'
$ws.Range("E4").Value = '''1.58%'
$ws.Range("G4").Value = '''23'
$ws.Range("D5").Value = '''0.05737'
$ws.Range("E5").Value = '''0.99%'
$ws.Range("G5").Value = '''23'
$ws.Range("D6").Value = '''6.564'
$ws.Range("E6").Value = '''1.18%'
$ws.Range("G6").Value = '''23'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.8587'
$ws.Range("E7").Value = '''4.71%'
$ws.Range("G7").Value = '''23'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '''0.8644'
$ws.Range("E8").Value = '''1.57%'
$ws.Range("G8").Value = '''23'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1367'
$ws.Range("E9").Value = '''2.59%'
$ws.Range("G9").Value = '''23'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07087'
$ws.Range("E10").Value = '''1.83%'
$ws.Range("G10").Value = '''23'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = '''0.03055'
$ws.Range("E11").Value = '''6.14%'
$ws.Range("G11").Value = '''23'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = '''0.09375'
$ws.Range("E12").Value = '''-0.15%'
$ws.Range("G12").Value = '''23'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = '''0.001542'
$ws.Range("E13").Value = '''2.21%'
$ws.Range("G13").Value = '''23'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0006035'
$ws.Range("E14").Value = '''0.57%'
$ws.Range("G14").Value = '''23'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006053'
$ws.Range("E15").Value = '''-2.61%'
$ws.Range("G15").Value = '''23'
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '''0.007489'
$ws.Range("E16").Value = '''5,224.85%'
$ws.Range("G16").Value = '''23'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.492'
$ws.Range("E17").Value = '''-0.59%'
$ws.Range("G17").Value = '''23'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''3.105'
$ws.Range("E18").Value = '''3.17%'
$ws.Range("G18").Value = '''23'
$ws.Range("E19").Value = '''-1.72%'
$ws.Range("G19").Value = '''23'
$ws.Range("E20").Value = '''1.34%'
$ws.Range("G20").Value = '''23'
$ws.Range("D21").Value = '''0.03331'
$ws.Range("E21").Value = '''2.86%'
$ws.Range("G21").Value = '''23'
$ws.Range("D22").Value = '''0.1291'
$ws.Range("E22").Value = '''1.32%'
$ws.Range("G22").Value = '''23'
$ws.Range("D23").Value = '''3.482'
$ws.Range("E23").Value = '''-2.01%'
$ws.Range("G23").Value = '''23'
$ws.Range("D24").Value = '''0.04153'
$ws.Range("E24").Value = '''2.90%'
$ws.Range("G24").Value = '''23'
$ws.Range("E25").Value = '''0.53%'
$ws.Range("G25").Value = '''23'
$ws.Range("D26").Value = '''0.001226'
$ws.Range("E26").Value = '''0.78%'
$ws.Range("G26").Value = '''23'
$ws.Range("D27").Value = '''0.004995'
$ws.Range("E27").Value = '''11.70%'
$ws.Range("G27").Value = '''23'
$ws.Range("E28").Value = '''2.56%'
$ws.Range("G28").Value = '''23'
$ws.Range("G29").Value = '''23'
$ws.Range("G30").Value = '''23'
$ws.Range("G31").Value = '''23'
$ws.Range("G32").Value = '''23'
$ws.Range("G33").Value = '''23'
$ws.Range("G34").Value = '''23'
$ws.Range("G35").Value = '''23'
$ws.Range("G36").Value = '''23'
$ws.Range("G37").Value = '''23'
$ws.Range("G38").Value = '''23'
$ws.Range("G39").Value = '''23'
$ws.Range("D40").Value = '''0.03754'
$ws.Range("E40").Value = '''0.99%'
$ws.Range("G40").Value = '''23'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1072'
$ws.Range("E41").Value = '''1.35%'
$ws.Range("G41").Value = '''23'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002451'
$ws.Range("E42").Value = '''-1.59%'
$ws.Range("G42").Value = '''23'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003522'
$ws.Range("E43").Value = '''-41.13%'
$ws.Range("G43").Value = '''23'
$ws.Range("D44").Value = '''0.009466'
$ws.Range("E44").Value = '''-1.34%'
$ws.Range("G44").Value = '''23'
$ws.Range("D45").Value = '''0.00005296'
$ws.Range("E45").Value = '''3.77%'
$ws.Range("G45").Value = '''23'
$ws.Range("G46").Value = '''23'
$ws.Range("D47").Value = '''0.05705'
$ws.Range("G47").Value = '''23'
$ws.Range("D48").Value = '''0.002281'
$ws.Range("E48").Value = '''-9.81%'
$ws.Range("G48").Value = '''23'
$ws.Range("G49").Value = '''23'
$ws.Range("G50").Value = '''23'
$ws.Range("G51").Value = '''23'
